$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price values are numeric-looking strings with a significant
# trailing zero (e.g. "0.8890"). Mark just those cells as Text first so Excel
# does not silently normalize them to "0.889" when the value is assigned.
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"

$ws.Range('D2').Value = '27.819.41'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '1.859.32'
$ws.Range('D4').Value = '1.037'
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').Value = '323.54'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').Value = '0.4415'
$ws.Range('E7').Value = '  +0.89%  '
$ws.Range('D8').Value = '0.3821'
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('D9').Value = '0.07447'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').Value = '0.8890'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').Value = '21.61'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = '1.861.26'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '5.543'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').Value = '6.750'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '0.07216'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '86.13'
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('D17').Value = '1.040'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = '0.000009121'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = '15.60'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = '27.830.83'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').Value = '5.312'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').Value = '11.28'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').Value = '2.098.30'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('D25').Value = '2.075'
$ws.Range('E25').Value = '  +6.58%  '
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').Value = '18.79'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').Value = '1.999'
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('D30').Value = '118.87'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').Value = '0.09127'
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.216'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '1.212'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').Value = '3.028'
$ws.Range('E34').Value = '  +5.23%  '
$ws.Range('D36').Value = '1.036'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').Value = '1.156'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '0.01988'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').Value = '0.05328'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('D40').Value = '2.859'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').Value = '0.5221'
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('D42').Value = '6.963'
$ws.Range('E42').Value = '  +3.35%  '
$ws.Range('D43').Value = '0.1681'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '8.806'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').Value = '110.88'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('D46').Value = '10.82'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').Value = '0.06579'
$ws.Range('D49').Value = '1.717'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').Value = '0.4741'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').Value = '1.886'
$ws.Range('E51').Value = '  +0.29%  '
